$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '69.234.89'
$ws.Range('E2').Value = '  -3.37%  '

Set-TextCell $ws.Range('D3') '3.516.62'
$ws.Range('E3').Value = '  -3.47%  '

Set-TextCell $ws.Range('D4') '0.998'
$ws.Range('E4').Value = '  -0.07%  '

Set-TextCell $ws.Range('D5') '577.54'
$ws.Range('E5').Value = '  -1.87%  '

Set-TextCell $ws.Range('D6') '171.11'
$ws.Range('E6').Value = '  -5.06%  '

Set-TextCell $ws.Range('D7') '0.617'
$ws.Range('E7').Value = '  +0.38%  '

Set-TextCell $ws.Range('D8') '3.507.81'
$ws.Range('E8').Value = '  -3.55%  '

$ws.Range('E9').Value = '  +0.15%  '

Set-TextCell $ws.Range('D10') '0.190'
$ws.Range('E10').Value = '  -6.40%  '

Set-TextCell $ws.Range('D11') '6.61'
$ws.Range('E11').Value = '  +12.48%  '

Set-TextCell $ws.Range('D12') '0.600'
$ws.Range('E12').Value = '  -1.23%  '

Set-TextCell $ws.Range('D13') '47.33'
$ws.Range('E13').Value = '  -4.86%  '

Set-TextCell $ws.Range('D14') '0.0000275'
$ws.Range('E14').Value = '  -3.66%  '

Set-TextCell $ws.Range('D15') '691.34'
$ws.Range('E15').Value = '  +1.08%  '

Set-TextCell $ws.Range('D16') '4.087.53'
$ws.Range('E16').Value = '  -3.29%  '

Set-TextCell $ws.Range('D17') '8.75'
$ws.Range('E17').Value = '  -2.72%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws.Range('D18') '69.108.94'
$ws.Range('E18').Value = '  -3.56%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range('D19') '3.513.59'
$ws.Range('E19').Value = '  -3.34%  '

$ws.Range('E20').Value = '  -1.59%  '

Set-TextCell $ws.Range('D21') '17.40'
$ws.Range('E21').Value = '  -5.02%  '

Set-TextCell $ws.Range('D22') '11.16'

Set-TextCell $ws.Range('D23') '0.913'
$ws.Range('E23').Value = '  -2.42%  '

Set-TextCell $ws.Range('D24') '16.60'
$ws.Range('E24').Value = '  -6.71%  '

Set-TextCell $ws.Range('D25') '97.66'
$ws.Range('E25').Value = '  -5.40%  '

Set-TextCell $ws.Range('D26') '3.83'
$ws.Range('E26').Value = '  -4.45%  '

$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range('D27') '2.68'
$ws.Range('E27').Value = '  -5.73%  '

$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range('D28') '1.00'
$ws.Range('E28').Value = '  -0.11%  '

Set-TextCell $ws.Range('D29') '9.43'
$ws.Range('E29').Value = '  -5.55%  '

Set-TextCell $ws.Range('D30') '33.15'
$ws.Range('E30').Value = '  -5.45%  '

Set-TextCell $ws.Range('D31') '8.86'
$ws.Range('E31').Value = '  -3.70%  '

Set-TextCell $ws.Range('D32') '3.18'
$ws.Range('E32').Value = '  -6.59%  '

Set-TextCell $ws.Range('D33') '7.30'
$ws.Range('E33').Value = '  +0.23%  '

Set-TextCell $ws.Range('D34') '1.34'
$ws.Range('E34').Value = '  -6.07%  '

Set-TextCell $ws.Range('D35') '3.77'
$ws.Range('E35').Value = '  -8.98%  '

Set-TextCell $ws.Range('D36') '567.95'
$ws.Range('E36').Value = '  -1.52%  '

Set-TextCell $ws.Range('D37') '10.88'
$ws.Range('E37').Value = '  -3.88%  '

$ws.Range('E38').Value = '  -3.68%  '

Set-TextCell $ws.Range('D39') '57.43'
$ws.Range('E39').Value = '  -3.41%  '

Set-TextCell $ws.Range('D40') '0.996'
$ws.Range('E40').Value = '  -0.16%  '

Set-TextCell $ws.Range('D41') '3.484.89'
$ws.Range('E41').Value = '  -5.11%  '

Set-TextCell $ws.Range('D42') '0.0441'
$ws.Range('E42').Value = '  -6.35%  '

$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range('D43') '0.138'
$ws.Range('E43').Value = '  -3.45%  '

$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws.Range('D44') '0.338'
$ws.Range('E44').Value = '  -2.52%  '

Set-TextCell $ws.Range('D45') '33.45'
$ws.Range('E45').Value = '  -6.25%  '

Set-TextCell $ws.Range('D46') '0.0₃0707'
$ws.Range('E46').Value = '  -7.09%  '

Set-TextCell $ws.Range('D47') '2.91'
$ws.Range('E47').Value = '  +4.01%  '

Set-TextCell $ws.Range('D48') '2.58'
$ws.Range('E48').Value = '  -6.33%  '

$ws.Range('E49').Value = '  -0.79%  '

Set-TextCell $ws.Range('D50') '133.84'
$ws.Range('E50').Value = '  +1.62%  '

$ws.Range('E51').Value = '  -0.54%  '
